# Stand doc list/detail views. Scoped all list views to stand/plot.

$wb = $excel.ActiveWorkbook
$props = $wb.Worksheets.Item("properties")

# The default view type changes from SPREADSHEET to LIST.
$props.Range("E2").Value = "LIST"

# Add the detail view filename property (row 3).
$props.Range("A3").Value = "Table"
$props.Range("B3").Value = "default"
$props.Range("C3").Value = "detailViewFileName"
$props.Range("D3").Value = "string"

# Add the list view filename property (row 4).
$props.Range("A4").Value = "Table"
$props.Range("B4").Value = "default"
$props.Range("C4").Value = "listViewFileName"
$props.Range("D4").Value = "string"

# Fill in the html paths last.
$props.Range("E3").Value = "config/tables/stand_doc/html/stand_doc_detail.html"
$props.Range("E4").Value = "config/tables/stand_doc/html/stand_doc_list.html"

# Select the new last cell and make this sheet the active one (matches the
# recorded author session ending on the properties sheet).
$props.Range("E4").Select() | Out-Null
$props.Activate() | Out-Null
